# Update status and tasks
#
# Marks a batch of robot-motor/sensor tasks complete (100%) with a
# completion date, adds a couple of notes about follow-up mechanical /
# sensor work, and leaves the cursor on C42 (matching where the author's
# Excel session ended up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18, 20-26: mark 100% complete and stamp the completion date.
$doneRows = 18,20,21,22,23,24,25,26
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 4).Value = 1        # column D - %Complete
    $ws.Cells.Item($r, 5).Value = 40062    # column E - Date Complete
}

# Notes (column F) for two of the newly-completed tasks.
$ws.Range("F24").Value = "May need to do some mechanical fixes for making sure the wheels rotate smoothly"
$ws.Range("F26").Value = "Use Bump Sensors"

# Leave the selection where the author left it.
$ws.Range("C42").Select()
